# Commit: "Don't strip two first rows of the sheet."
#
# Sheets "E" and "I" each had two superfluous leading blank rows (rows 1-2)
# sitting above their real data (which started at row 3). This removes
# those two leading blank rows so the real data shifts up to start at
# row 1, and moves the active workbook tab/selection from sheet "E" to
# the last sheet "Z" (with its selection reset to A1).

$wb = $excel.ActiveWorkbook

# --- Sheet "E": drop the two leading blank rows, data moves up to row 1 ---
$wsE = $wb.Worksheets.Item("E")
$wsE.Rows("1:2").Delete()
[void]$wsE.Range("A1").Select()

# --- Sheet "I": drop the two leading blank rows, data moves up to row 1 ---
$wsI = $wb.Worksheets.Item("I")
$wsI.Rows("1:2").Delete()
[void]$wsI.Range("A1").Select()

# --- Sheet "Z" becomes the active / selected tab, selection reset to A1 ---
$wsZ = $wb.Worksheets.Item("Z")
$wsZ.Activate()
[void]$wsZ.Range("A1").Select()
